# Updated cryptos list on Sat Mar 23 22:11:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''64.999.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +3.44%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.389.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +3.20%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.18%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''559.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +4.00%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''174.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +2.95%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +2.46%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.378.94'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +3.12%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -0.02%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +10.98%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  +4.71%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''54.02'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +3.89%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +6.20%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''9.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +3.93%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.930.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +3.09%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''18.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +2.33%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  +3.10%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.377.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +2.62%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''64.910.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +3.25%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''11.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +2.67%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +3.33%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''469.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +14.03%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +12.17%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''4.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +3.86%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''86.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +5.65%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''13.51'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +2.33%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''2.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +8.92%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''10.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +2.92%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  +3.39%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''30.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +6.92%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''6.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +7.98%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''11.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +2.26%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''571.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -0.23%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''61.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +7.00%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  +2.98%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +0.09%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''3.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +5.97%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -3.32%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''35.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +2.42%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.0₃0743'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.31%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +2.85%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''3.088.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +0.06%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -0.14%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +4.59%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.0416'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +4.99%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  +5.91%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  +3.10%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''3.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -2.49%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''2.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +0.91%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''139.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +5.25%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''8.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +4.14%  '
$ws.Range("E51").Style = "Normal"
